$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 3: Version
$ws.Range("B3").Value = "0.1.7"

# Row 6: Status
$ws.Range("B6").Value = "draft"

# Row 8: Date
$ws.Range("B8").Value = "2024-08-23T10:17:11-05:00"

# Row 10: Contact (organization) - update text
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Row 11: Contact (person) - update text (was duplicate "No display for ContactDetail")
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Insert a brand new row 12 for Jurisdiction, pushing Description/Purpose/Copyright/Immutable down
$ws.Rows("12:12").Insert()

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

$ws.Range("A12:B12").WrapText = $true
$ws.Range("A12:B12").VerticalAlignment = -4160
$ws.Range("A12:B12").Borders.Item(9).LineStyle = 1
$ws.Range("A12:B12").Borders.Item(9).Weight = 2
